$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains text-formatted numbers (e.g. "1.003"); force text type
# so Excel does not auto-convert numeric-looking strings to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.854.13'
$ws.Range('E2').Value = '  -2.78%  '
$ws.Range('D3').Value = '1.792.26'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').Value = '315.88'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').Value = '0.5398'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.3781'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.07411'
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('D10').Value = '41.85'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('D11').Value = '1.084'
$ws.Range('E11').Value = '  -2.86%  '
$ws.Range('D12').Value = '1.004'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').Value = '6.181'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').Value = '7.368'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '20.22'
$ws.Range('E15').Value = '  -3.51%  '
$ws.Range('D16').Value = '1.798.47'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = '88.67'
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '0.00001058'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').Value = '0.06496'
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').Value = '17.27'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = '5.915'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').Value = '27.929.21'
$ws.Range('E23').Value = '  -2.50%  '
$ws.Range('D24').Value = '11.09'
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('D25').Value = '2.093'
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('D26').Value = '155.74'
$ws.Range('E26').Value = '  -3.05%  '
$ws.Range('D27').Value = '20.31'
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('D28').Value = '2.004.15'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('D29').Value = '2.328'
$ws.Range('E29').Value = '  -2.05%  '
$ws.Range('D30').Value = '121.71'
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('D31').Value = '0.1103'
$ws.Range('E31').Value = '  +5.65%  '
$ws.Range('D32').Value = '1.103'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('D33').Value = '3.665'
$ws.Range('E33').Value = '  -0.78%  '
$ws.Range('D34').Value = '5.508'
$ws.Range('E34').Value = '  -2.87%  '
$ws.Range('D35').Value = '0.06919'
$ws.Range('E35').Value = '  +6.40%  '
$ws.Range('D36').Value = '0.2198'
$ws.Range('E36').Value = '  -2.79%  '
$ws.Range('D37').Value = '0.02275'
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('D38').Value = '5.039'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('D39').Value = '8.418'
$ws.Range('E39').Value = '  -5.31%  '
$ws.Range('D40').Value = '11.20'
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').Value = '0.6098'
$ws.Range('E41').Value = '  -2.53%  '
$ws.Range('D44').Value = '13.28'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('D45').Value = '3.678'
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('D46').Value = '0.5688'
$ws.Range('E46').Value = '  -3.26%  '
$ws.Range('D47').Value = '124.10'
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').Value = '1.174'
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('D49').Value = '1.904'
$ws.Range('E49').Value = '  -2.96%  '
$ws.Range('D50').Value = '0.06799'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('B42').Value = 'WEMIXTOKEN'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').Value = '1.420'
$ws.Range('E42').Value = '  +1.80%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '1.162'
$ws.Range('E43').Value = '  -4.06%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '71.53'
$ws.Range('E51').Value = '  -1.65%  '

# Restore original (default) cell style on column D now that values are text,
# so no stray number-format style id is left on these cells.
$dRange.Style = "Normal"
